$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final updates for VA creation:
# - the vendor id in A2 changes from 800577 to 800021
# - the two trailing vendor rows (A3: 801740, A4: 803260) are removed,
#   shrinking the used range down to A1:A2
$ws.Range("A3:A4").ClearContents()
$ws.Range("A2").Value = 800021
